# Updates the "cryptos" worksheet with refreshed price/volume data,
# including a couple of ranking swaps (rows 25/26 and rows 43/45).
#
# Price values in column D are stored as plain TEXT in the workbook
# (e.g. "6.50", "0.0000160", "1.00") even though they look numeric.
# Writing such a string straight into a Range.Value makes Excel
# auto-coerce it into a real number (dropping trailing zeros /
# switching to scientific notation), so for any new value that would
# parse as a number we force text by prefixing an apostrophe and then
# reset the cell style back to the default ("Normal") so no stray
# number-format gets baked in.

function Set-TextValue($sheet, $cellRef, $val) {
    $sheet.Range($cellRef).Value = "'" + $val
    $sheet.Range($cellRef).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '64.741.40'
$ws.Range("E2").Value = '  -1.79%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '3.112.95'
$ws.Range("E3").Value = '  -7.78%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  -0.01%  '

# Row 5 - BNB
Set-TextValue $ws 'D5' '565.22'
$ws.Range("E5").Value = '  -2.72%  '

# Row 6 - Solana
Set-TextValue $ws 'D6' '167.48'
$ws.Range("E6").Value = '  -6.45%  '

# Row 8 - XRP
Set-TextValue $ws 'D8' '0.593'
$ws.Range("E8").Value = '  -5.08%  '

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = '3.112.15'
$ws.Range("E9").Value = '  -7.74%  '

# Row 10 - Dogecoin
Set-TextValue $ws 'D10' '0.122'
$ws.Range("E10").Value = '  -6.70%  '

# Row 11 - Toncoin
Set-TextValue $ws 'D11' '6.50'
$ws.Range("E11").Value = '  -6.21%  '

# Row 12 - Cardano
Set-TextValue $ws 'D12' '0.385'
$ws.Range("E12").Value = '  -6.79%  '

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = '3.640.94'
$ws.Range("E13").Value = '  -8.11%  '

# Row 14 - TRON
$ws.Range("E14").Value = '  +0.48%  '

# Row 15 - Avalanche
Set-TextValue $ws 'D15' '26.62'
$ws.Range("E15").Value = '  -8.05%  '

# Row 16 - WrappedBTC
$ws.Range("D16").Value = '64.698.90'
$ws.Range("E16").Value = '  -2.01%  '

# Row 17 - ShibaInu
Set-TextValue $ws 'D17' '0.0000160'
$ws.Range("E17").Value = '  -6.47%  '

# Row 18 - WrappedEther
$ws.Range("D18").Value = '3.108.83'
$ws.Range("E18").Value = '  -7.66%  '

# Row 19 - Polkadot
Set-TextValue $ws 'D19' '5.60'
$ws.Range("E19").Value = '  -4.03%  '

# Row 20 - Chainlink
Set-TextValue $ws 'D20' '12.60'
$ws.Range("E20").Value = '  -7.45%  '

# Row 21 - BitcoinCash
Set-TextValue $ws 'D21' '351.73'
$ws.Range("E21").Value = '  -4.07%  '

# Row 22 - Uniswap
Set-TextValue $ws 'D22' '7.12'
$ws.Range("E22").Value = '  -5.31%  '

# Row 23 - Dai
Set-TextValue $ws 'D23' '1.00'
$ws.Range("E23").Value = '  +0.37%  '

# Row 24 - Litecoin
Set-TextValue $ws 'D24' '68.15'
$ws.Range("E24").Value = '  -6.20%  '

# Row 25 - was WrappedeETH, now Polygon (ranking swap with row 26)
$ws.Range("B25").Value = 'Polygon'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws 'D25' '0.486'
$ws.Range("E25").Value = '  -7.96%  '

# Row 26 - was Polygon, now WrappedeETH (ranking swap with row 25)
$ws.Range("B26").Value = 'WrappedeETH'
$ws.Range("C26").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D26").Value = '3.249.50'
$ws.Range("E26").Value = '  -7.83%  '

# Row 27 - PEPE
Set-TextValue $ws 'D27' '0.0000110'
$ws.Range("E27").Value = '  -10.82%  '

# Row 28 - InternetComputer(DFINITY)
Set-TextValue $ws 'D28' '9.48'
$ws.Range("E28").Value = '  -2.37%  '

# Row 29 - Kaspa
$ws.Range("E29").Value = '  -2.63%  '

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = '  -0.02%  '

# Row 31 - USDe
$ws.Range("E31").Value = '  -0.18%  '

# Row 32 - PancakeSwap
Set-TextValue $ws 'D32' '1.89'
$ws.Range("E32").Value = '  -4.46%  '

# Row 33 - EthereumClassic
Set-TextValue $ws 'D33' '21.53'
$ws.Range("E33").Value = '  -6.65%  '

# Row 34 - NEARProtocol
Set-TextValue $ws 'D34' '5.19'
$ws.Range("E34").Value = '  -9.58%  '

# Row 35 - Aptos
Set-TextValue $ws 'D35' '6.49'
$ws.Range("E35").Value = '  -7.16%  '

# Row 36 - Fetch.AI
Set-TextValue $ws 'D36' '1.17'
$ws.Range("E36").Value = '  -5.45%  '

# Row 37 - Monero
Set-TextValue $ws 'D37' '158.16'
$ws.Range("E37").Value = '  -1.82%  '

# Row 38 - ImmutableX
Set-TextValue $ws 'D38' '1.41'
$ws.Range("E38").Value = '  -7.17%  '

# Row 39 - Mantle
Set-TextValue $ws 'D39' '0.814'
$ws.Range("E39").Value = '  -4.77%  '

# Row 40 - EnergySwap
Set-TextValue $ws 'D40' '25.71'
$ws.Range("E40").Value = '  -5.16%  '

# Row 41 - Stacks
Set-TextValue $ws 'D41' '1.72'
$ws.Range("E41").Value = '  -3.01%  '

# Row 42 - Maker
$ws.Range("D42").Value = '2.617.35'
$ws.Range("E42").Value = '  -2.47%  '

# Row 43 - was OKB, now dogwifhat (ranking swap with row 45)
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws 'D43' '2.38'
$ws.Range("E43").Value = '  -8.38%  '

# Row 44 - RenderToken
Set-TextValue $ws 'D44' '5.98'
$ws.Range("E44").Value = '  -5.37%  '

# Row 45 - was dogwifhat, now OKB (ranking swap with row 43)
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws 'D45' '39.22'
$ws.Range("E45").Value = '  -1.26%  '

# Row 46 - Filecoin
Set-TextValue $ws 'D46' '4.10'
$ws.Range("E46").Value = '  -5.21%  '

# Row 47 - Hedera
Set-TextValue $ws 'D47' '0.0645'
$ws.Range("E47").Value = '  -4.16%  '

# Row 48 - InjectiveProtocol
Set-TextValue $ws 'D48' '23.39'
$ws.Range("E48").Value = '  -3.95%  '

# Row 49 - Bittensor
Set-TextValue $ws 'D49' '313.85'
$ws.Range("E49").Value = '  -6.43%  '

# Row 50 - VeChain
Set-TextValue $ws 'D50' '0.0267'
$ws.Range("E50").Value = '  -5.38%  '

# Row 51 - Stellar
Set-TextValue $ws 'D51' '0.101'
$ws.Range("E51").Value = '  -3.16%  '
